$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 0.6579926714760105
$ws.Range("D6").Value = 0.009068124681414892
$ws.Range("E6").Value = 0.1928860392336463
$ws.Range("F6").Value = 0.1928303336857762
$ws.Range("G6").Value = 0.007163578701580939
$ws.Range("H6").Value = 0.647948717948718
$ws.Range("K6").Value = 0.2972766307864244
$ws.Range("L6").Value = 0.2968197459314154
$ws.Range("M6").Value = 0.01086837734671845
$ws.Range("N6").Value = 0.6536018402839749
$ws.Range("P6").Value = 0.01354743174035911

$ws.Range("B7").Value = 0.6509588193708024
$ws.Range("C7").Value = 0.6526159886550158
$ws.Range("D7").Value = 0.007746899931491852
$ws.Range("E7").Value = 0.1833035044166041
$ws.Range("F7").Value = 0.1840305737424639
$ws.Range("G7").Value = 0.004991028202449576
$ws.Range("H7").Value = 0.6576222038111019
$ws.Range("K7").Value = 0.286694658901079
$ws.Range("L7").Value = 0.2874993469672502
$ws.Range("M7").Value = 0.007480595057322747
$ws.Range("N7").Value = 0.6538927123516881
$ws.Range("O7").Value = 0.6545485626997557
$ws.Range("P7").Value = 0.009643825894607794

$ws.Range("B8").Value = 0.6527483209614705
$ws.Range("C8").Value = 0.6541403842711692
$ws.Range("D8").Value = 0.007706481349815865
$ws.Range("E8").Value = 0.1843112244897959
$ws.Range("F8").Value = 0.1847661948929004
$ws.Range("G8").Value = 0.004766417463952714
$ws.Range("K8").Value = 0.2880050733828592
$ws.Range("L8").Value = 0.2883934120039368
$ws.Range("M8").Value = 0.007162110424629433
$ws.Range("N8").Value = 0.6552590844304127
$ws.Range("O8").Value = 0.6554015632549682
$ws.Range("P8").Value = 0.009307095007076292

$ws.Range("B9").Value = 0.6415893736402153
$ws.Range("C9").Value = 0.6386404287443808
$ws.Range("D9").Value = 0.00877482679191423
$ws.Range("E9").Value = 0.1885838150289017
$ws.Range("F9").Value = 0.1868204389404011
$ws.Range("G9").Value = 0.00529533616469915
$ws.Range("H9").Value = 0.6692307692307692
$ws.Range("I9").Value = 0.6669230769230768
$ws.Range("J9").Value = 0.02288237078699806
$ws.Range("K9").Value = 0.294250281848929
$ws.Range("L9").Value = 0.291843190867129
$ws.Range("M9").Value = 0.008074280600331423
$ws.Range("N9").Value = 0.6536731314573541
$ws.Range("O9").Value = 0.6510045482587272
$ws.Range("P9").Value = 0.01062930004420306

$ws.Range("B10").Value = 0.6574201305393336
$ws.Range("C10").Value = 0.6583930747671711
$ws.Range("D10").Value = 0.008791020344452588
$ws.Range("E10").Value = 0.1924513915364087
$ws.Range("F10").Value = 0.192852410911758
$ws.Range("G10").Value = 0.006938299989315343
$ws.Range("I10").Value = 0.6464102564102564
$ws.Range("J10").Value = 0.02545315791534024
$ws.Range("K10").Value = 0.296679400528945
$ws.Range("L10").Value = 0.2970403587972042
$ws.Range("M10").Value = 0.01046870287223622
$ws.Range("N10").Value = 0.6529433140976064
$ws.Range("O10").Value = 0.6531546835630062
$ws.Range("P10").Value = 0.01298307511292331

$ws.Range("C12").Value = 0.6515542790040783
$ws.Range("D12").Value = 0.008348164708974783
$ws.Range("F12").Value = 0.149961287847496
$ws.Range("G12").Value = 0.004493427459083526
$ws.Range("L12").Value = 0.2406645366855023
$ws.Range("M12").Value = 0.007149463416398572
$ws.Range("O12").Value = 0.6326080643066874
$ws.Range("P12").Value = 0.01016789111197173

$ws.Range("B17").Value = 0.6015831134564644
$ws.Range("E17").Value = 0.1561938958707361
$ws.Range("K17").Value = 0.2531277276694792
$ws.Range("N17").Value = 0.6306879036965547

$ws.Range("B21").Value = 0.6535937738440409
$ws.Range("E21").Value = 0.1708502024291498
$ws.Range("K21").Value = 0.271034039820167
$ws.Range("N21").Value = 0.6543447777839678

$ws.Range("B22").Value = 0.6595319397018405
$ws.Range("C22").Value = 0.6547830646472935
$ws.Range("D22").Value = 0.008316448542932573
$ws.Range("E22").Value = 0.2020224719101124
$ws.Range("F22").Value = 0.1996340740909303
$ws.Range("G22").Value = 0.007067378239773102
$ws.Range("H22").Value = 0.6260445682451253
$ws.Range("K22").Value = 0.3054706082229018
$ws.Range("L22").Value = 0.3027795804075785
$ws.Range("M22").Value = 0.01051749000355718
$ws.Range("N22").Value = 0.6450627722504598
$ws.Range("O22").Value = 0.6426761442266942
$ws.Range("P22").Value = 0.01246509873020882

$ws.Range("B23").Value = 0.6834680005130178
$ws.Range("E23").Value = 0.2099447513812155
$ws.Range("K23").Value = 0.3159645232815964
$ws.Range("N23").Value = 0.6638042084761224

$ws.Range("B24").Value = 0.6646897810218978
$ws.Range("C24").Value = 0.6653728203766472
$ws.Range("D24").Value = 0.009681022000855258
$ws.Range("E24").Value = 0.182896379525593
$ws.Range("F24").Value = 0.1799600160647329
$ws.Range("G24").Value = 0.01431046136318765
$ws.Range("K24").Value = 0.2850194552529183
$ws.Range("L24").Value = 0.2797049762502483
$ws.Range("M24").Value = 0.02241595691260664
$ws.Range("N24").Value = 0.6561477844660412
$ws.Range("O24").Value = 0.648657364687899
$ws.Range("P24").Value = 0.02776317239831637

$ws.Range("B26").Value = 0.6264456658651093
$ws.Range("C26").Value = 0.6307397718825986
$ws.Range("D26").Value = 0.00797425324314611
$ws.Range("E26").Value = 0.1814636863728153
$ws.Range("F26").Value = 0.1811430265738616
$ws.Range("G26").Value = 0.004114947950153339
$ws.Range("I26").Value = 0.6553846153846153
$ws.Range("J26").Value = 0.02260776661041757
$ws.Range("K26").Value = 0.2854170089261267
$ws.Range("L26").Value = 0.2838001624225816
$ws.Range("M26").Value = 0.006541186530918769
$ws.Range("N26").Value = 0.644701300890396
$ws.Range("O26").Value = 0.6415135830490373
$ws.Range("P26").Value = 0.009028001722758468

$ws.Range("B27").Value = 0.5980762624527654
$ws.Range("C27").Value = 0.6009960053825372
$ws.Range("D27").Value = 0.013296488608787
$ws.Range("E27").Value = 0.1771934292627022
$ws.Range("F27").Value = 0.1762726431973381
$ws.Range("G27").Value = 0.004115429118951222
$ws.Range("H27").Value = 0.7135897435897436
$ws.Range("I27").Value = 0.7002564102564103
$ws.Range("J27").Value = 0.02386682313481048
$ws.Range("K27").Value = 0.2838926859124757
$ws.Range("L27").Value = 0.2815787756193535
$ws.Range("M27").Value = 0.005842553675592371
$ws.Range("N27").Value = 0.6485743252622603
$ws.Range("O27").Value = 0.6443888771076597
$ws.Range("P27").Value = 0.008279314525923885
